# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for the rows
# whose dialog-act classification changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 9;   Tag = "ba"; Label = "Appreciation" },
    @{ Row = 10;  Tag = "b";  Label = "Acknowledge (Backchannel)" },
    @{ Row = 14;  Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 49;  Tag = "b";  Label = "Acknowledge (Backchannel)" },
    @{ Row = 80;  Tag = "b";  Label = "Acknowledge (Backchannel)" },
    @{ Row = 91;  Tag = "%";  Label = "Uninterpretable" },
    @{ Row = 101; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 107; Tag = "%";  Label = "Uninterpretable" },
    @{ Row = 112; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 115; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 116; Tag = "b";  Label = "Acknowledge (Backchannel)" },
    @{ Row = 133; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 167; Tag = "aa"; Label = "Agree/Accept" },
    @{ Row = 170; Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 174; Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 175; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 179; Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 205; Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 207; Tag = "sv"; Label = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Label
}
